$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'5.37%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'34.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'12.79%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.190"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'4.75%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07824"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'6.21%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.313"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-1.12%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'8.054"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'4.24%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.989"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'7.10%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9266"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1.69%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1001"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'9.57%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1835"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.08547"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'3.55%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03368"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'8.04%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09918"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.47%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001482"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-1.18%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.04648"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'2.93%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.005737"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.55%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.486"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.24%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'2.104"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.43%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.3441"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'3.46%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1324"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'3.23%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.584"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'10.27%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.2385"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'13.74%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'0.97%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'6.41%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001301"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'0.11%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003402"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'0.29%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01748"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'10.78%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04748"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'6.24%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007700"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'4.44%"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'6.43%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007079"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-25.67%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002210"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-1.26%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.01004"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'15.38%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00005988"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-1.72%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.12%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'3.880"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'51.26%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.002691"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'34.67%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.12%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.12%"
$ws.Range("E51").Style = "Normal"
